# Updated Test cases content Multiple specimen, derivative, aliquot.
# Targets the "MySQL: https://...caTissue Database Dump/v2.0/MySQL and deploy
# application." paragraph (Pre-requisites section).

$d = $word.ActiveDocument

# Locate the target paragraph precisely (scope all Find operations to it so
# the similarly-worded "Oracle:" paragraph right above is never touched).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "MySQL*and deploy application.*") {
        $target = $p
        break
    }
}

$pr = $target.Range

# 1) Paragraph formatting: add adjustRightInd=0 and switch spacing to an
#    explicit single-line rule (w:line="240" w:lineRule="auto").
$pr.ParagraphFormat.AutoAdjustRightIndent = $false
$pr.ParagraphFormat.LineSpacingRule = 0

# 2) Split "...Dump/v2.0/" into two runs: "...Dump/v2" + ".0/" by nudging a
#    character property on just the ".0/" tail (forces the engine to break
#    the run there) and then restoring it so there is no visible formatting
#    change.
$pr2 = $target.Range
$found = $pr2.Find.Execute(".0/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pr2.Font.Bold = 1
$pr2.Font.Bold = 0

# 3) Expand the trailing sentence with the Label Generator setting detail.
$pr3 = $target.Range
$found3 = $pr3.Find.Execute(" and deploy application.", $true, $false, $false, $false, $false, $true, 1, $false, " and deploy application with Label Generator Setting ON for Specimen “edu.wustl.catissuecore.namegenerator.DefaultSpecimenLabelGenerator“", 2)

Write-Host "edit complete"
